# Updated cryptos list on Fri Mar 22 11:51:55 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.595.64'
$ws.Range("E2").Value = '  -3.74%  '
$ws.Range("D3").Value = '3.424.06'
$ws.Range("E3").Value = '  -3.06%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = '''570.61'
$ws.Range("E5").Value = '  +2.06%  '
$ws.Range("D6").Value = '''173.57'
$ws.Range("D7").Value = '''0.619'
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D9").Value = '''0.623'
$ws.Range("E9").Value = '  -1.45%  '
$ws.Range("D10").Value = '''0.156'
$ws.Range("E10").Value = '  +2.91%  '
$ws.Range("D11").Value = '''55.17'
$ws.Range("E11").Value = '  +0.55%  '
$ws.Range("D12").Value = '''0.0000270'
$ws.Range("E12").Value = '  -0.10%  '
$ws.Range("D13").Value = '''9.10'
$ws.Range("E13").Value = '  -3.10%  '
$ws.Range("D14").Value = '3.979.92'
$ws.Range("E14").Value = '  -2.77%  '
$ws.Range("E15").Value = '  -1.13%  '
$ws.Range("D16").Value = '3.431.28'
$ws.Range("E16").Value = '  -2.93%  '
$ws.Range("D17").Value = '''18.02'
$ws.Range("E17").Value = '  -1.25%  '
$ws.Range("D18").Value = '''11.87'
$ws.Range("E18").Value = '  -1.70%  '
$ws.Range("D19").Value = '64.710.13'
$ws.Range("E19").Value = '  -3.56%  '
$ws.Range("D20").Value = '''0.988'
$ws.Range("E20").Value = '  -0.88%  '
$ws.Range("D21").Value = '''406.64'
$ws.Range("E21").Value = '  -5.56%  '
$ws.Range("D22").Value = '''4.17'
$ws.Range("E22").Value = '  +1.83%  '
$ws.Range("D23").Value = '''4.31'
$ws.Range("E23").Value = '  +4.40%  '
$ws.Range("D24").Value = '''83.18'
$ws.Range("E24").Value = '  -2.56%  '
$ws.Range("D25").Value = '''13.10'
$ws.Range("E25").Value = '  +6.14%  '
$ws.Range("D26").Value = '''10.80'
$ws.Range("E26").Value = '  -2.75%  '
$ws.Range("D27").Value = '''2.78'
$ws.Range("E27").Value = '  -3.99%  '
$ws.Range("D28").Value = '''5.99'
$ws.Range("E28").Value = '  -2.28%  '
$ws.Range("D29").Value = '''8.90'
$ws.Range("E29").Value = '  -2.27%  '
$ws.Range("D30").Value = '''29.74'
$ws.Range("E30").Value = '  -2.39%  '
$ws.Range("D31").Value = '''6.70'
$ws.Range("E31").Value = '  +1.94%  '
$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D32").Value = '''586.01'
$ws.Range("E32").Value = '  -8.62%  '
$ws.Range("B33").Value = 'Cosmos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D33").Value = '''11.47'
$ws.Range("E33").Value = '  -2.28%  '
$ws.Range("E34").Value = '  -2.82%  '
$ws.Range("D35").Value = '''59.09'
$ws.Range("E35").Value = '  -1.59%  '
$ws.Range("E36").Value = '  +3.55%  '
$ws.Range("D37").Value = '''1.00'
$ws.Range("E37").Value = '  +0.08%  '
$ws.Range("D38").Value = '''36.01'
$ws.Range("E38").Value = '  -6.32%  '
$ws.Range("D39").Value = '0.0₃0757'
$ws.Range("E39").Value = '  -6.47%  '
$ws.Range("D40").Value = '''3.46'
$ws.Range("E40").Value = '  +2.64%  '
$ws.Range("D41").Value = '''0.375'
$ws.Range("E41").Value = '  -3.68%  '
$ws.Range("D42").Value = '3.175.97'
$ws.Range("E42").Value = '  +4.80%  '
$ws.Range("D43").Value = '''0.999'
$ws.Range("E43").Value = '  -0.01%  '
$ws.Range("D44").Value = '''2.90'
$ws.Range("E44").Value = '  +0.71%  '
$ws.Range("D45").Value = '''2.49'
$ws.Range("E45").Value = '  -5.89%  '
$ws.Range("D46").Value = '''3.20'
$ws.Range("E46").Value = '  -4.17%  '
$ws.Range("D47").Value = '''0.0407'
$ws.Range("E47").Value = '  -2.80%  '
$ws.Range("E48").Value = '  -4.48%  '
$ws.Range("E49").Value = '  -1.66%  '
$ws.Range("E50").Value = '  -2.74%  '
$ws.Range("D51").Value = '''136.28'
$ws.Range("E51").Value = '  -5.01%  '
